# Actualizar 02-05-2021 09-48-34
# Refresh pass: nudge the previous batch's timestamp by a hair (re-save
# jitter) and append a brand-new batch of 14 rows (one per monitored
# service) with the new check timestamp + hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. tiny float nudge on the previous batch's Fecha column (rows 478-491) ---
$prevTimestamp = 44232.3876172338
for ($r = 478; $r -le 491; $r++) {
    $ws.Cells.Item($r, 4).Value = $prevTimestamp
}

# --- 2. append the new batch (rows 492-505) ---
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Target address actually stored in the relationship (the fragment after
# "#" becomes the hyperlink's SubAddress/location instead of being part
# of the stored Target).
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkSubAddresses = @("","","","","","","","","/","","","","","")

$newTimestamp = 44232.40867991291
$startRow = 492

for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $urls[$i]
    $ws.Cells.Item($row, 3).Value = "Disponible"
    $ws.Cells.Item($row, 4).Value = $newTimestamp

    $cell = $ws.Cells.Item($row, 2)
    if ($linkSubAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($cell, $linkAddresses[$i], $linkSubAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($cell, $linkAddresses[$i])
    }
}

# Re-apply the workbook's existing "Hyperlink" cell look (font/underline,
# no border) to the freshly linked B cells so they match the style already
# used by every other hyperlink cell in the sheet instead of keeping
# whatever default formatting Hyperlinks.Add left behind.
$ws.Range("B2:B2").Copy() | Out-Null
$ws.Range("B492:B505").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-apply the existing "Fecha" date-time number format to the new D cells
# (brand-new rows start out with General formatting otherwise).
$ws.Range("D2:D2").Copy() | Out-Null
$ws.Range("D492:D505").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
